# Applies the edits described by the diff:
#  - "Base cases and interventions" -> "Cases and interventions"
#  - Case 1(a)/1(b)/2/3/4 paragraphs renumbered to Case 1/2/3/4/5
#    (each now refers back to the previous case by number, e.g.
#    "as Case 1(a)" -> "as Case 1,")
#  - Table 1 and Table 2 case-label cells renumbered the same way
#    (1a/1b/2/3/4 -> 1/2/3/4/5)
#
# (The diff also contains several hunks where Word's automatic
# spelling/grammar checker re-flowed already-identical text into extra
# <w:r> runs wrapped in <w:proofErr .../> markers, e.g. "(water volume
# change)", "(morphological change)" and "...maintenance dreding." --
# those carry no visible text change and are not achievable through the
# documented Word object model, so they are left alone.)

$d = $word.ActiveDocument

$script:failures = 0

function Replace-InRange($range, $findText, $replaceText, $label) {
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    # Replace:=1 (wdReplaceOne) -- replaces only the first match found
    # inside $range. (wdReplaceAll ignores the supplied range scope in
    # this runtime and would rewrite every matching occurrence in the
    # whole document, which is not what we want here since several of
    # these strings - e.g. the "1a - ..." row labels - are duplicated
    # verbatim between Table 1 and Table 2.)
    $result = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 0, $false, $replaceText, 1)
    if (-not $result) {
        Write-Output "WARNING: replacement failed for $label (find=[$findText])"
        $script:failures = $script:failures + 1
    }
    return $result
}

$dash = [string][char]0x2013   # en dash "-" used throughout the document

# ---------------------------------------------------------------
# 1. Heading
# ---------------------------------------------------------------
$rng = $d.Paragraphs(2).Range
$find = "Base cases and interventions"
$repl = "Cases and interventions"
Replace-InRange $rng $find $repl "heading" | Out-Null

# ---------------------------------------------------------------
# 2. "Base cases and interventions" paragraphs (Case 1(a) .. Case 4)
# ---------------------------------------------------------------

# Case 1(a) -> Case 1
$rng = $d.Paragraphs(3).Range
$find = "Case 1(a) " + $dash + " sea level rise of 2mm/year"
$repl = "Case 1" + $dash + " sea level rise of 2mm/year"
Replace-InRange $rng $find $repl "case1" | Out-Null

# Case 1(b) -> Case 2
$rng = $d.Paragraphs(4).Range
$find = "Case 1(b) " + $dash + " sea level rise with a nodal cycle of amplitude 0.15m"
$repl = "Case 2 " + $dash + " sea level rise with a nodal cycle of amplitude 0.15m"
Replace-InRange $rng $find $repl "case1b" | Out-Null

# Case 2 -> Case 3 (as Case 1(a) -> as Case 1,)
$rng = $d.Paragraphs(5).Range
$find = "Case 2 " + $dash + " as Case 1(a) with the historic interventions included (but no maintenance dredging)"
$repl = "Case 3 " + $dash + " as Case 1, with the historic interventions included (but no maintenance dredging)"
Replace-InRange $rng $find $repl "case2" | Out-Null

# Case 3 -> Case 4 (as Case 2 -> as Case 3,)
$rng = $d.Paragraphs(6).Range
$find = "Case 3 " + $dash + " as Case 2 with a capital dredge of 7.2"
$repl = "Case 4 " + $dash + " as Case 3,with a capital dredge of 7.2"
Replace-InRange $rng $find $repl "case3" | Out-Null

# Case 4 -> Case 5 (as Case 3 -> as Case 4,)
$rng = $d.Paragraphs(7).Range
$find = "Case 4 " + $dash + " as Case 3 with reclamation of the inner flats removing 20 Ha with a volume of 0.5 Mm"
$repl = "Case 5 " + $dash + " as Case 4, with reclamation of the inner flats removing 20 Ha with a volume of 0.5 Mm"
Replace-InRange $rng $find $repl "case4" | Out-Null

# ---------------------------------------------------------------
# 3. Table case-label cells, for both Table 1 (water level rise only)
#    and Table 2 (incl. maintenance dredging) -- identical relabelling.
# ---------------------------------------------------------------
foreach ($tableIndex in 2, 4) {
    $t = $d.Tables($tableIndex)

    $rng = $t.Cell(2,1).Range
    $find = "1a " + $dash + " sea level rise of 2mm/year"
    $repl = "1 " + $dash + " sea level rise of 2mm/year"
    Replace-InRange $rng $find $repl "table$tableIndex-1a" | Out-Null

    $rng = $t.Cell(3,1).Range
    $find = "1b " + $dash + " As 1a with ntc of 0.15m"
    $repl = "2 " + $dash + " As 1, with ntc of 0.15m"
    Replace-InRange $rng $find $repl "table$tableIndex-1b" | Out-Null

    $rng = $t.Cell(4,1).Range
    $find = "2 " + $dash + " historic changes + slr"
    $repl = "3 " + $dash + " historic changes + slr"
    Replace-InRange $rng $find $repl "table$tableIndex-2" | Out-Null

    $rng = $t.Cell(5,1).Range
    $find = "3 " + $dash + " As 2 with dredge in 2000*"
    $repl = "4 " + $dash + " As 3, with dredge in 2000*"
    Replace-InRange $rng $find $repl "table$tableIndex-3" | Out-Null

    $rng = $t.Cell(6,1).Range
    $find = "4 " + $dash + " As 3 with reclamation in 2020*"
    $repl = "5 " + $dash + " As 4, with reclamation in 2020*"
    Replace-InRange $rng $find $repl "table$tableIndex-4" | Out-Null
}

if ($script:failures -eq 0) {
    Write-Output "Done: all replacements applied successfully"
} else {
    Write-Output "Done with $script:failures failure(s) - see warnings above"
}
